# Add new "Genus" feature column (J) with header value 8, and refresh the
# evaluation metric columns (E:I) with re-run values, plus the new J column
# of metric values for each classification row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell J1 = 8, styled like the rest of row 1 (bold/centered) ---
$ws.Range("J1").Value = 8
$ws.Range("J1").Style = $ws.Range("I1").Style

# --- Updated evaluation metrics for existing columns E:I, plus new column J ---

# Row 3 (Phylum)
$ws.Range("E3").Value = 0.6422762488378444
$ws.Range("F3").Value = 0.7592434738368422
$ws.Range("G3").Value = 0.7512697025727029
$ws.Range("H3").Value = 0.9562995693712001
$ws.Range("I3").Value = 0.9606374937356758
$ws.Range("J3").Value = 0.9591339735637363

# Row 4 (Class)
$ws.Range("E4").Value = 0.5641166287940037
$ws.Range("F4").Value = 0.7067203585025023
$ws.Range("G4").Value = 0.6964537641577606
$ws.Range("H4").Value = 0.9451823546881468
$ws.Range("I4").Value = 0.9471612920030504
$ws.Range("J4").Value = 0.9481911627870443

# Row 5 (Order)
$ws.Range("E5").Value = 0.3307204436357077
$ws.Range("F5").Value = 0.432783673130987
$ws.Range("G5").Value = 0.3670866409380765
$ws.Range("H5").Value = 0.683352271985583
$ws.Range("I5").Value = 0.6845876244994668
$ws.Range("J5").Value = 0.6682234216999524

# Row 6 (Family)
$ws.Range("E6").Value = 0.3658359288155295
$ws.Range("F6").Value = 0.5238323792101567
$ws.Range("G6").Value = 0.4544159848709955
$ws.Range("H6").Value = 0.8516830987559882
$ws.Range("I6").Value = 0.8530601670637019
$ws.Range("J6").Value = 0.8495801154964431

# Row 7 (Genus)
$ws.Range("E7").Value = 0.3262781025224333
$ws.Range("F7").Value = 0.4801323453186883
$ws.Range("G7").Value = 0.4293958368657118
$ws.Range("H7").Value = 0.9065692010125441
$ws.Range("I7").Value = 0.9047803169031409
$ws.Range("J7").Value = 0.9039027077099593
